$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.168.18'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '1.838.37'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.53'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6302'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07504'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2930'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.15'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07752'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '1.833.45'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.997'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6698'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.96%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '82.68'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009365'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.017'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '29.164.33'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = '2.077.91'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.60'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '223.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.007'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.144'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.006'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '160.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1400'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.518'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.98'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.506'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05923'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +13.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.159'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.071'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.204'
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7505'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.851'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.141'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.686'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').Value = '1.229.06'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.769'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01795'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.563'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8951'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.008'
$ws.Range('D43').Style = "Normal"
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '102.34'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('B45').Value = 'XinFinNetwork'
$ws.Range('C45').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.08035'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +18.62%  '
$ws.Range('D46').Value = '1.978.59'
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000124'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '65.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.5112'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4074'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.017'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.82%  '
